$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 515:516, pushing the existing weekly Coliflor
# records (old rows 515-548) down to 517-550.
$ws.Rows("515:516").Insert()

# Row 515 - Primera (new weekly record)
$ws.Range("A515").Value = 3
$ws.Range("B515").Value = "Femacal de La Calera"
$ws.Range("C515").Value = "Coquimbo"
$ws.Range("D515").Value = 44610
$ws.Range("E515").Value = 5
$ws.Range("F515").Value = 100112008
$ws.Range("G515").Value = "Coliflor"
$ws.Range("H515").Value = "Sin especificar"
$ws.Range("I515").Value = "Primera"
$ws.Range("J515").Value = 550
$ws.Range("K515").Value = 1200
$ws.Range("L515").Value = 1200
$ws.Range("M515").Value = 1200
$ws.Range("N515").Value = "$/unidad"
$ws.Range("O515").Value = "Provincia de Quillota"
$ws.Range("P515").Value = 1200
$ws.Range("Q515").Value = 1
$ws.Range("R515").Value = "Hortaliza"

# Row 516 - Segunda (new weekly record)
$ws.Range("A516").Value = 3
$ws.Range("B516").Value = "Femacal de La Calera"
$ws.Range("C516").Value = "Coquimbo"
$ws.Range("D516").Value = 44610
$ws.Range("E516").Value = 5
$ws.Range("F516").Value = 100112008
$ws.Range("G516").Value = "Coliflor"
$ws.Range("H516").Value = "Sin especificar"
$ws.Range("I516").Value = "Segunda"
$ws.Range("J516").Value = 750
$ws.Range("K516").Value = 800
$ws.Range("L516").Value = 800
$ws.Range("M516").Value = 800
$ws.Range("N516").Value = "$/unidad"
$ws.Range("O516").Value = "Provincia de Quillota"
$ws.Range("P516").Value = 800
$ws.Range("Q516").Value = 1
$ws.Range("R516").Value = "Hortaliza"
